$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# This paragraph (in div p036v_3) originally reads (display text, tags shown
# literally as escaped text):
#   ... longue &amp; grosse <tl>verge de <m>fer</m></tl> ronde, puys la retirent<lb/> ...
# and must become:
#   ... longue &amp; <tl>grosse verge de <m>fer</m> ronde</tl>, puys la retirent<lb/> ...
# i.e. "grosse " moves from before <tl> to just after it (joining "verge de "),
# and the closing </tl> moves from immediately after </m> to after " ronde"
# (the <tl> markup span widens to also cover "grosse" and "ronde").
# ---------------------------------------------------------------------------

# 1) " grosse " -> " "  (the word "grosse" is removed from before the <tl> tag)
$r1 = $d.Content
$r1.Find.Execute(" grosse ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r1.Text = " "

# 2) "verge de " -> "grosse verge de "  (the word "grosse" re-appears, now after <tl>)
$r2 = $d.Content
$r2.Find.Execute("verge de ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$r2.Text = "grosse verge de "

# 3) "</m></tl>" -> "</m>"  (drop the trailing </tl> from this run; it gets
#    reinserted after " ronde" in step 4)
$r3 = $d.Content
$r3.Find.Execute("fer</m></tl> ronde", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tagOnly = $d.Range($r3.Start + 3, $r3.Start + 12)
$tagOnly.Text = "</m>"

# 4) " ronde, puys la retirent" -> " ronde</tl>, puys la retirent"
#    (insert a new, separately-formatted "</tl>" run right after " ronde")
$r4 = $d.Content
$r4.Find.Execute(" ronde, puys la retirent", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $r4.Start + 6   # right after " ronde", before the comma

# grab formatted text (text + full run formatting) of an existing closing
# "</tl>" run elsewhere in the document to use as the template for the new run
$tagRef = $d.Content
$tagRef.Find.Execute("<tl>four</tl>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tagRefClose = $d.Range($tagRef.End - 5, $tagRef.End)
$tagFormatted = $tagRefClose.FormattedText

$insertionPoint = $d.Range($splitPoint, $splitPoint)
$insertionPoint.FormattedText = $tagFormatted

Write-Output "done"
